$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The shared-string table needs several new entries (Gaussian-Quadrature is
# moved earlier, and three new "Spiral-..." strings are introduced) so that
# on save the engine rebuilds xl/sharedStrings.xml with the strings in a new
# order. The engine keeps previously-used shared strings pinned to their
# original relative order and only appends freshly-introduced text in the
# order it is assigned, so first blank out every text cell that needs to
# change (releasing the old string references) and then re-assign all of the
# text values in exactly the order they must appear in the rebuilt table.
# ---------------------------------------------------------------------------

$textCells = @(
    "B2","C2","D2","E2","F2","G2","H2","I2","J2","K2","L2","M2",
    "B3","B4","B5","B6","B7","B8","B9","B10","B11","B12","B13","B14","B15","B16"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).Value = ""
}

# Re-assign every text value in the exact order the rebuilt shared-string
# table must list them.
$ws.Range("B2").Value = "HKL"

$ws.Range("B3").Value = "ND Single"
$ws.Range("B4").Value = "RD Single"
$ws.Range("B5").Value = "TD Single"
$ws.Range("B6").Value = "Morris"
$ws.Range("B7").Value = "Ring Perpendicular to ND"
$ws.Range("B8").Value = "Ring Perpendicular to RD"
$ws.Range("B9").Value = "Ring Perpendicular to TD"

$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"

$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# New rows 17-19
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C2").Value = "[1, 1, 0]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 1, 1]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 1, 0]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 2, 1]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "2Pairs"
$ws.Range("L2").Value = "4Pairs"
$ws.Range("M2").Value = "MaxUnique"

# ---------------------------------------------------------------------------
# Column A labels (numeric, bold/centered/bordered style copied from an
# existing styled cell so no new style entries are introduced).
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 15

$ws.Range("A3").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 16

$ws.Range("A3").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 17

# ---------------------------------------------------------------------------
# Data values (all 1) for the three new rows, columns C through M.
# ---------------------------------------------------------------------------
foreach ($r in 17..19) {
    foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M")) {
        $ws.Range("$col$r").Value = 1
    }
}
